$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 30..65 in column A get replaced with incentive-vars:1000 .. incentive-vars:1035
for ($row = 30; $row -le 65; $row++) {
    $id = 1000 + ($row - 30)
    $ws.Cells.Item($row, 1).Value = "incentive-vars:$id"
}
